$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-44: bump the "changed" date from 2023-09-03
# to 2023-09-06 for every record on the sheet.
$newDate = Get-Date -Year 2023 -Month 9 -Day 6 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("C2:C44").Value = $newDate
